$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Set the value for H7 (Identifie code smells, Day 4 -> 0.4 effort logged)
$ws.Range("H7").Value = 0.4

# Update the active selection to match the authored change
$ws.Range("H8").Select()

$excel.Calculate()

$wb.Save()
